$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits right after the paragraph that ends
# with "...ne" vektorski procesor podacima ". It needs to be moved to the
# author/co-author line, right after "Nikola Kovačević, " once the
# co-author "Vuk Vranjković" is removed. Remove it from its old spot first
# so the text edit below doesn't have to worry about it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Remove the co-author "Vuk Vranjković, " from the author line, keeping the
# same run formatting for the remaining "Nikola Kovačević, " text.
$rng = $d.Content
$rng.Find.Execute("Nikola Kovačević, Vuk Vranjković, ", $true, $false, $false, $false, $false, $true, 1, $false, "Nikola Kovačević, ", 2)

# Re-insert the "_GoBack" bookmark right after the (now shorter) author
# text, before the "Fakultet tehničkih nauka, Novi Sad" run.
$insertPoint = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $insertPoint)
